$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "1__211006132800_Waves_001.txt"

$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 513.08
$ws.Range("G2").Value = 7.71
$ws.Range("H2").Value = 5.78
$ws.Range("I2").Value = 0.79
$ws.Range("J2").Value = 0.06
$ws.Range("K2").Value = 3.46
$ws.Range("L2").Value = 0.71
$ws.Range("M2").Value = 0.05
$ws.Range("N2").Value = 3.86
$ws.Range("O2").Value = 0.55
$ws.Range("P2").Value = 0.04
$ws.Range("Q2").Value = 44.63
$ws.Range("R2").Value = 15.41
$ws.Range("S2").Value = 1.15
$ws.Range("T2").Value = 2.75
$ws.Range("U2").Value = 0.6
$ws.Range("V2").Value = 0.04
$ws.Range("W2").Value = 132.64
$ws.Range("X2").Value = 26.53
$ws.Range("Y2").Value = 1.97
$ws.Range("Z2").Value = 8.35
$ws.Range("AA2").Value = 1.16
$ws.Range("AB2").Value = 0.09
$ws.Range("AC2").Value = 7.27
$ws.Range("AD2").Value = 1.03
$ws.Range("AE2").Value = 0.08
$ws.Range("AF2").Value = 11.36
$ws.Range("AG2").Value = 2.3
$ws.Range("AH2").Value = 0.17
$ws.Range("AI2").Value = 17.43
$ws.Range("AJ2").Value = 0.99
$ws.Range("AK2").Value = 0.07000000000000001
